# Auto-generated edit script applying numeric corrections to the Seraph_Profits workbook.
# Each sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) holds a Leve-profit table in columns A:N.
# H:N are recomputed price/profit columns; this script rewrites the updated figures cell-by-cell
# and clears the handful of cells that became blank (no HQ data) after the refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 177.33333   # H4
$ws.Cells.Item(4, 9).Value = 177.33333   # I4
$ws.Cells.Item(4, 11).Value = 177.33333   # K4
$ws.Cells.Item(4, 13).Value = -63.33332999999999   # M4
$ws.Cells.Item(6, 8).Value = 250.25   # H6
$ws.Cells.Item(6, 9).Value = 250.25   # I6
$ws.Cells.Item(6, 11).Value = 750.75   # K6
$ws.Cells.Item(6, 13).Value = -638.75   # M6
$ws.Cells.Item(17, 8).Value = 1031.5   # H17
$ws.Cells.Item(17, 9).Value = 1031.5   # I17
$ws.Cells.Item(17, 11).Value = 3094.5   # K17
$ws.Cells.Item(17, 13).Value = -2926.5   # M17
$ws.Cells.Item(39, 8).Value = 10.4   # H39
$ws.Cells.Item(39, 10).Value = 0   # J39
$ws.Cells.Item(39, 12).Value = 0   # L39
$ws.Cells.Item(39, 14).ClearContents()   # N39
$ws.Cells.Item(41, 8).Value = 199.33333   # H41
$ws.Cells.Item(41, 9).Value = 246.16667   # I41
$ws.Cells.Item(41, 11).Value = 246.16667   # K41
$ws.Cells.Item(41, 13).Value = 193.83333   # M41
$ws.Cells.Item(44, 8).Value = 20050   # H44
$ws.Cells.Item(44, 10).Value = 20050   # J44
$ws.Cells.Item(44, 12).Value = 20050   # L44
$ws.Cells.Item(44, 14).Value = -20974   # N44
$ws.Cells.Item(53, 8).Value = 112.47059   # H53
$ws.Cells.Item(53, 9).Value = 85   # I53
$ws.Cells.Item(53, 10).Value = 178.4   # J53
$ws.Cells.Item(53, 11).Value = 85   # K53
$ws.Cells.Item(53, 12).Value = 178.4   # L53
$ws.Cells.Item(53, 13).Value = 552   # M53
$ws.Cells.Item(53, 14).Value = -1452.4   # N53
$ws.Cells.Item(132, 8).Value = 1880.5   # H132
$ws.Cells.Item(132, 9).Value = 1072.8667   # I132
$ws.Cells.Item(132, 11).Value = 3218.6001   # K132
$ws.Cells.Item(132, 13).Value = -688.6001000000001   # M132

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 1990.6154   # H26
$ws.Cells.Item(26, 9).Value = 1237.8   # I26
$ws.Cells.Item(26, 10).Value = 4500   # J26
$ws.Cells.Item(26, 11).Value = 1237.8   # K26
$ws.Cells.Item(26, 12).Value = 4500   # L26
$ws.Cells.Item(26, 13).Value = -907.8   # M26
$ws.Cells.Item(26, 14).Value = -5160   # N26
$ws.Cells.Item(32, 8).Value = 3678.0256   # H32
$ws.Cells.Item(32, 9).Value = 2384.1428   # I32
$ws.Cells.Item(32, 11).Value = 2384.1428   # K32
$ws.Cells.Item(32, 13).Value = -2097.1428   # M32
$ws.Cells.Item(46, 8).Value = 15833.333   # H46
$ws.Cells.Item(46, 9).Value = 17750   # I46
$ws.Cells.Item(46, 10).Value = 12000   # J46
$ws.Cells.Item(46, 11).Value = 17750   # K46
$ws.Cells.Item(46, 12).Value = 12000   # L46
$ws.Cells.Item(46, 13).Value = -17431   # M46
$ws.Cells.Item(46, 14).Value = -12638   # N46
$ws.Cells.Item(61, 8).Value = 7427.5   # H61
$ws.Cells.Item(61, 9).Value = 7697.222   # I61
$ws.Cells.Item(61, 11).Value = 7697.222   # K61
$ws.Cells.Item(61, 13).Value = -7485.222   # M61
$ws.Cells.Item(132, 8).Value = 4393.0625   # H132
$ws.Cells.Item(132, 9).Value = 4614.5   # I132
$ws.Cells.Item(132, 11).Value = 13843.5   # K132
$ws.Cells.Item(132, 13).Value = -11313.5   # M132
$ws.Cells.Item(136, 8).Value = 7427.5   # H136
$ws.Cells.Item(136, 9).Value = 7697.222   # I136
$ws.Cells.Item(136, 11).Value = 23091.666   # K136
$ws.Cells.Item(136, 13).Value = -20541.666   # M136

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3495   # H20
$ws.Cells.Item(20, 9).Value = 3495   # I20
$ws.Cells.Item(20, 11).Value = 3495   # K20
$ws.Cells.Item(20, 13).Value = -3248   # M20
$ws.Cells.Item(105, 8).Value = 3116.074   # H105
$ws.Cells.Item(105, 9).Value = 2549.3684   # I105
$ws.Cells.Item(105, 11).Value = 2549.3684   # K105
$ws.Cells.Item(105, 13).Value = -802.3683999999998   # M105
$ws.Cells.Item(134, 8).Value = 0   # H134
$ws.Cells.Item(134, 9).Value = 0   # I134
$ws.Cells.Item(134, 11).Value = 0   # K134
$ws.Cells.Item(134, 13).ClearContents()   # M134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(29, 8).Value = 28333   # H29
$ws.Cells.Item(29, 9).Value = 28000   # I29
$ws.Cells.Item(29, 10).Value = 28499.5   # J29
$ws.Cells.Item(29, 11).Value = 28000   # K29
$ws.Cells.Item(29, 12).Value = 28499.5   # L29
$ws.Cells.Item(29, 13).Value = -27707   # M29
$ws.Cells.Item(29, 14).Value = -29085.5   # N29
$ws.Cells.Item(58, 8).Value = 2599.7273   # H58
$ws.Cells.Item(58, 9).Value = 1845.2222   # I58
$ws.Cells.Item(58, 11).Value = 1845.2222   # K58
$ws.Cells.Item(58, 13).Value = -1642.2222   # M58
$ws.Cells.Item(106, 8).Value = 12500   # H106
$ws.Cells.Item(106, 10).Value = 12500   # J106
$ws.Cells.Item(106, 12).Value = 12500   # L106
$ws.Cells.Item(106, 14).Value = -15024   # N106
$ws.Cells.Item(136, 8).Value = 2599.7273   # H136
$ws.Cells.Item(136, 9).Value = 1845.2222   # I136
$ws.Cells.Item(136, 11).Value = 5535.6666   # K136
$ws.Cells.Item(136, 13).Value = -2985.6666   # M136

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 336.91666   # H2
$ws.Cells.Item(2, 10).Value = 542.5714   # J2
$ws.Cells.Item(2, 12).Value = 3255.4284   # L2
$ws.Cells.Item(2, 14).Value = -3481.4284   # N2
$ws.Cells.Item(11, 8).Value = 1499.6666   # H11
$ws.Cells.Item(11, 9).Value = 1250   # I11
$ws.Cells.Item(11, 11).Value = 3750   # K11
$ws.Cells.Item(11, 13).Value = -3610   # M11
$ws.Cells.Item(108, 8).Value = 9342.857   # H108
$ws.Cells.Item(108, 9).Value = 800   # I108
$ws.Cells.Item(108, 11).Value = 2400   # K108
$ws.Cells.Item(108, 13).Value = 480   # M108

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 6500   # H7
$ws.Cells.Item(7, 9).Value = 0   # I7
$ws.Cells.Item(7, 10).Value = 6500   # J7
$ws.Cells.Item(7, 11).Value = 0   # K7
$ws.Cells.Item(7, 12).Value = 6500   # L7
$ws.Cells.Item(7, 13).ClearContents()   # M7
$ws.Cells.Item(7, 14).Value = -6724   # N7
$ws.Cells.Item(8, 8).Value = 6500   # H8
$ws.Cells.Item(8, 9).Value = 0   # I8
$ws.Cells.Item(8, 10).Value = 6500   # J8
$ws.Cells.Item(8, 11).Value = 0   # K8
$ws.Cells.Item(8, 12).Value = 6500   # L8
$ws.Cells.Item(8, 13).ClearContents()   # M8
$ws.Cells.Item(8, 14).Value = -6778   # N8
$ws.Cells.Item(39, 8).Value = 55000   # H39
$ws.Cells.Item(39, 10).Value = 55000   # J39
$ws.Cells.Item(39, 12).Value = 55000   # L39
$ws.Cells.Item(39, 14).Value = -56064   # N39
$ws.Cells.Item(102, 8).Value = 1129.1428   # H102
$ws.Cells.Item(102, 9).Value = 982.63635   # I102
$ws.Cells.Item(102, 11).Value = 982.63635   # K102
$ws.Cells.Item(102, 13).Value = 639.36365   # M102
$ws.Cells.Item(113, 8).Value = 11958.714   # H113
$ws.Cells.Item(113, 9).Value = 903.6667   # I113
$ws.Cells.Item(113, 11).Value = 903.6667   # K113
$ws.Cells.Item(113, 13).Value = 1266.3333   # M113
$ws.Cells.Item(126, 8).Value = 5434.7   # H126
$ws.Cells.Item(126, 9).Value = 4921.3335   # I126
$ws.Cells.Item(126, 10).Value = 6204.75   # J126
$ws.Cells.Item(126, 11).Value = 14764.0005   # K126
$ws.Cells.Item(126, 12).Value = 18614.25   # L126
$ws.Cells.Item(126, 13).Value = -12294.0005   # M126
$ws.Cells.Item(126, 14).Value = -23554.25   # N126
$ws.Cells.Item(132, 8).Value = 2890.1853   # H132
$ws.Cells.Item(132, 9).Value = 3123.261   # I132
$ws.Cells.Item(132, 11).Value = 9369.782999999999   # K132
$ws.Cells.Item(132, 13).Value = -6839.782999999999   # M132
$ws.Cells.Item(141, 8).Value = 39998.5   # H141
$ws.Cells.Item(141, 10).Value = 39998.5   # J141
$ws.Cells.Item(141, 12).Value = 39998.5   # L141
$ws.Cells.Item(141, 14).Value = -50358.5   # N141

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value = 4375.143   # H3
$ws.Cells.Item(3, 10).Value = 4375.143   # J3
$ws.Cells.Item(3, 12).Value = 4375.143   # L3
$ws.Cells.Item(3, 14).Value = -4599.143   # N3
$ws.Cells.Item(7, 8).Value = 3201.2   # H7
$ws.Cells.Item(7, 9).Value = 3999   # I7
$ws.Cells.Item(7, 11).Value = 3999   # K7
$ws.Cells.Item(7, 13).Value = -3887   # M7
$ws.Cells.Item(12, 8).Value = 0   # H12
$ws.Cells.Item(12, 10).Value = 0   # J12
$ws.Cells.Item(12, 12).Value = 0   # L12
$ws.Cells.Item(12, 14).ClearContents()   # N12
$ws.Cells.Item(15, 8).Value = 4375.143   # H15
$ws.Cells.Item(15, 10).Value = 4375.143   # J15
$ws.Cells.Item(15, 12).Value = 4375.143   # L15
$ws.Cells.Item(15, 14).Value = -4715.143   # N15
$ws.Cells.Item(40, 8).Value = 3214.6316   # H40
$ws.Cells.Item(40, 9).Value = 3214.6316   # I40
$ws.Cells.Item(40, 11).Value = 3214.6316   # K40
$ws.Cells.Item(40, 13).Value = -3078.6316   # M40
$ws.Cells.Item(80, 8).Value = 32000   # H80
$ws.Cells.Item(80, 10).Value = 32000   # J80
$ws.Cells.Item(80, 12).Value = 32000   # L80
$ws.Cells.Item(80, 14).Value = -34246   # N80
$ws.Cells.Item(83, 8).Value = 32000   # H83
$ws.Cells.Item(83, 10).Value = 32000   # J83
$ws.Cells.Item(83, 12).Value = 96000   # L83
$ws.Cells.Item(83, 14).Value = -107232   # N83
$ws.Cells.Item(100, 8).Value = 2509.2   # H100
$ws.Cells.Item(100, 9).Value = 2349   # I100
$ws.Cells.Item(100, 10).Value = 3150   # J100
$ws.Cells.Item(100, 11).Value = 2349   # K100
$ws.Cells.Item(100, 12).Value = 3150   # L100
$ws.Cells.Item(100, 13).Value = -1808   # M100
$ws.Cells.Item(100, 14).Value = -4232   # N100
$ws.Cells.Item(126, 8).Value = 3201.2   # H126
$ws.Cells.Item(126, 9).Value = 3999   # I126
$ws.Cells.Item(126, 11).Value = 11997   # K126
$ws.Cells.Item(126, 13).Value = -9527   # M126
$ws.Cells.Item(128, 8).Value = 0   # H128
$ws.Cells.Item(128, 10).Value = 0   # J128
$ws.Cells.Item(128, 12).Value = 0   # L128
$ws.Cells.Item(128, 14).ClearContents()   # N128
$ws.Cells.Item(132, 8).Value = 0   # H132
$ws.Cells.Item(132, 9).Value = 0   # I132
$ws.Cells.Item(132, 10).Value = 0   # J132
$ws.Cells.Item(132, 11).Value = 0   # K132
$ws.Cells.Item(132, 12).Value = 0   # L132
$ws.Cells.Item(132, 13).ClearContents()   # M132
$ws.Cells.Item(132, 14).ClearContents()   # N132
$ws.Cells.Item(136, 8).Value = 7145.273   # H136
$ws.Cells.Item(136, 9).Value = 7109.8   # I136
$ws.Cells.Item(136, 11).Value = 21329.4   # K136
$ws.Cells.Item(136, 13).Value = -18779.4   # M136

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(10, 8).Value = 100000   # H10
$ws.Cells.Item(10, 9).Value = 100000   # I10
$ws.Cells.Item(10, 11).Value = 100000   # K10
$ws.Cells.Item(10, 13).Value = -99831   # M10
$ws.Cells.Item(11, 8).Value = 16999.5   # H11
$ws.Cells.Item(11, 10).Value = 16999.5   # J11
$ws.Cells.Item(11, 12).Value = 16999.5   # L11
$ws.Cells.Item(11, 14).Value = -17283.5   # N11
$ws.Cells.Item(13, 8).Value = 650   # H13
$ws.Cells.Item(13, 9).Value = 650   # I13
$ws.Cells.Item(13, 11).Value = 650   # K13
$ws.Cells.Item(13, 13).Value = -510   # M13
$ws.Cells.Item(14, 8).Value = 16998   # H14
$ws.Cells.Item(14, 10).Value = 16998   # J14
$ws.Cells.Item(14, 12).Value = 16998   # L14
$ws.Cells.Item(14, 14).Value = -17334   # N14
$ws.Cells.Item(25, 8).Value = 21438.5   # H25
$ws.Cells.Item(25, 10).Value = 21438.5   # J25
$ws.Cells.Item(25, 12).Value = 21438.5   # L25
$ws.Cells.Item(25, 14).Value = -22024.5   # N25
$ws.Cells.Item(93, 8).Value = 44999.5   # H93
$ws.Cells.Item(93, 10).Value = 44999.5   # J93
$ws.Cells.Item(93, 12).Value = 44999.5   # L93
$ws.Cells.Item(93, 14).Value = -49991.5   # N93
$ws.Cells.Item(100, 8).Value = 2314.2856   # H100
$ws.Cells.Item(100, 10).Value = 2376.5   # J100
$ws.Cells.Item(100, 12).Value = 4753   # L100
$ws.Cells.Item(100, 14).Value = -5835   # N100
$ws.Cells.Item(132, 8).Value = 2575.75   # H132
$ws.Cells.Item(132, 9).Value = 2575.75   # I132
$ws.Cells.Item(132, 10).Value = 0   # J132
$ws.Cells.Item(132, 11).Value = 7727.25   # K132
$ws.Cells.Item(132, 12).Value = 0   # L132
$ws.Cells.Item(132, 13).Value = -5197.25   # M132
$ws.Cells.Item(132, 14).ClearContents()   # N132
